# Update the "dSF" column (column F) values in Sheet1 to reflect the
# repulled / recalculated data from the commit "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = -4
    3  = -2
    4  = -3
    5  = 0
    6  = 1
    7  = 4
    8  = -6
    9  = 2
    10 = 3
    12 = -4
    13 = -6
    14 = 7
    15 = 3
    16 = -2
    17 = 1
    18 = -2
    20 = 2
    21 = -1
    22 = 4
    23 = -1
    24 = -3
    25 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
